$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two PRODUCT values in column B throughout the data rows.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -eq "Tropical Punch Workout Boost") {
        $cell.Value2 = "LIV_Tropical Punch Workout Boost"
    } elseif ($val -eq "Hawaiian Punch Pre-Workout") {
        $cell.Value2 = "C4_Hawaiian Punch Pre-Workout"
    }
}

# Update the sheet view / selection state.
$ws.Range("J10").Select()

Write-Output "done"
